$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("potential action") before the existing explanation column,
# which shifts the old column B (explanation) to column C.
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "potential action"

$ws.Range("B2").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B3").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B4").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B5").Value = "no action needed"
$ws.Range("B6").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B7").Value = "re-scrape all calls from alternative source, drop company if can't find CNHI"
$ws.Range("B8").Value = "add as many calls as possible from alternative source, no action needed on failure"
$ws.Range("B9").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B10").Value = "drop company"
$ws.Range("B11").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B12").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B13").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B14").Value = "add as many calls as possible from alternative source, no action needed on failure"
$ws.Range("B15").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B16").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B17").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B18").Value = "add as many calls as possible from alternative source, no action needed on failure"
$ws.Range("B19").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B20").Value = "add as many calls as possible from alternative source, no action needed on failure"
$ws.Range("B21").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B22").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B23").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B24").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B25").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B26").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B27").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B28").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B29").Value = "add as many calls as possible from alternative source, no action needed on failure"
$ws.Range("B30").Value = "add as many calls as possible from alternative source, no action needed on failure"
$ws.Range("B31").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B32").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B33").Value = "add as many calls as possible from alternative source, no action needed on failure"
$ws.Range("B34").Value = "scrape all calls for this ticker with no date restriction"
$ws.Range("B35").Value = "re-scrape all calls from alternative source, drop company if can't find quarterly"

# Update the CNHI explanation text (row 7) to describe the Raven Industries situation.
$ws.Range("C7").Value = "The calls are for Raven Industries, which was acquired by CNHI in 2021. We can drop or try to get correct CNHI calls"

# Resize columns B and C (closest achievable widths under this engine's column-width quantization).
$ws.Columns.Item(2).ColumnWidth = 66.0
$ws.Columns.Item(3).ColumnWidth = 234.45

# Apply an AutoFilter across the full data range and register the hidden _FilterDatabase name.
$ws.Range("A1:C35").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$35")
$filterName.Visible = $false

# Move the active selection to match the saved workbook view.
$ws.Range("A28").Select()
